$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------
# Status text update: "In Translation" -> "Handed back: in sync with en-US"
# (Status column = C on both language sheets; the Overview sheet also
# mirrors the same status string in its "zh-cn"/"de-de" columns E/F,
# since they shared the same string in the workbook.)
# ---------------------------------------------------------------
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# zh-cn sheet: fill in "Latest Target File" (I) and "Latest Handback File" (J)
# for both data rows, and add the hyperlink on I (same target as the
# corresponding "Source File Name" hyperlink in column A).
# ---------------------------------------------------------------
$zhcn.Range("I2").Value = "a48795bd-0d8f-4434-b32b-2ec949890ad7.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a21239bae42f8bf42a9e1d6aced818b81d9187c2/e2e/a48795bd-0d8f-4434-b32b-2ec949890ad7.md", "", "", "a48795bd-0d8f-4434-b32b-2ec949890ad7.md") | Out-Null
$zhcn.Range("J2").Value = "a48795bd-0d8f-4434-b32b-2ec949890ad7.c124d74ee52885663b162c028aafa7037b6283da.zh-cn.xlf"

$zhcn.Range("I3").Value = "a804360b-9b49-41c4-9e2c-5eb2ceb089e4.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a21239bae42f8bf42a9e1d6aced818b81d9187c2/e2e/a804360b-9b49-41c4-9e2c-5eb2ceb089e4.md", "", "", "a804360b-9b49-41c4-9e2c-5eb2ceb089e4.md") | Out-Null
$zhcn.Range("J3").Value = "a804360b-9b49-41c4-9e2c-5eb2ceb089e4.2307d90c699b6506a3a69b2466e4c94d9408452b.zh-cn.xlf"

# "Latest Handback DateTime" (K) had the never-handed-back placeholder; update it now.
$zhcn.Range("K2").Value = "2016-08-12 08:28:40"
$zhcn.Range("K3").Value = "2016-08-12 08:28:40"

# ---------------------------------------------------------------
# de-de sheet: same shape of update, but a distinct handback timestamp.
# ---------------------------------------------------------------
$dede.Range("I2").Value = "a48795bd-0d8f-4434-b32b-2ec949890ad7.md"
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a21239bae42f8bf42a9e1d6aced818b81d9187c2/e2e/a48795bd-0d8f-4434-b32b-2ec949890ad7.md", "", "", "a48795bd-0d8f-4434-b32b-2ec949890ad7.md") | Out-Null
$dede.Range("J2").Value = "a48795bd-0d8f-4434-b32b-2ec949890ad7.c124d74ee52885663b162c028aafa7037b6283da.de-de.xlf"

$dede.Range("I3").Value = "a804360b-9b49-41c4-9e2c-5eb2ceb089e4.md"
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a21239bae42f8bf42a9e1d6aced818b81d9187c2/e2e/a804360b-9b49-41c4-9e2c-5eb2ceb089e4.md", "", "", "a804360b-9b49-41c4-9e2c-5eb2ceb089e4.md") | Out-Null
$dede.Range("J3").Value = "a804360b-9b49-41c4-9e2c-5eb2ceb089e4.2307d90c699b6506a3a69b2466e4c94d9408452b.de-de.xlf"

$dede.Range("K2").Value = "2016-08-12 08:28:49"
$dede.Range("K3").Value = "2016-08-12 08:28:49"

# ---------------------------------------------------------------
# Widen the columns that now hold file names / long hyperlink text
# so the handback report is readable.
# ---------------------------------------------------------------
$overview.Range("E1").ColumnWidth = 29.9777047293527
$overview.Range("F1").ColumnWidth = 29.9777047293527

$zhcn.Range("C1").ColumnWidth = 29.9777047293527
$zhcn.Range("I1").ColumnWidth = 40
$zhcn.Range("J1").ColumnWidth = 40

$dede.Range("C1").ColumnWidth = 29.9777047293527
$dede.Range("I1").ColumnWidth = 40
$dede.Range("J1").ColumnWidth = 40

$wb.Save()
